# Appends the "PRUEBAS DE MANTENIMIENTO" section at the end of the
# document, right after the final "Caso de uso: ..." paragraph and
# before the section break (w:sectPr).
#
# New paragraphs added (all carrying lang="es-ES", matching the rest
# of the document):
#   1. (empty)
#   2. "PRUEBAS DE MANTENIMIENTO:"
#   3. "Se realizan posterior a la entrega del producto"
#   4. "Es realizada cuando el sistema sistema sufre cambios "
#   5. (empty)

$d = $word.ActiveDocument

function Get-DocEndRange($doc) {
    $count = $doc.Paragraphs.Count
    $para = $doc.Paragraphs.Item($count)
    $r = $para.Range
    $r.Collapse(0)   # wdCollapseEnd
    return $r
}

function Add-TailParagraph($doc, [string]$text) {
    $r = Get-DocEndRange $doc
    # A leading paragraph mark ends the current last paragraph and
    # starts a brand new one; anything after it becomes that new
    # paragraph's text.
    $r.Text = "`r" + $text
    if ($text.Length -gt 0) {
        $newPara = $doc.Paragraphs.Item($doc.Paragraphs.Count)
        $newPara.Range.LanguageID = "es-ES"
    }
}

Add-TailParagraph $d ""
Add-TailParagraph $d "PRUEBAS DE MANTENIMIENTO:"
Add-TailParagraph $d "Se realizan posterior a la entrega del producto"
Add-TailParagraph $d "Es realizada cuando el sistema sistema sufre cambios "
Add-TailParagraph $d ""
